$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 132 (shifts existing rows 132-170 down to 133-171)
$ws.Rows.Item(132).Insert()

# Populate the newly inserted row 132 with the new weekly price entry
$ws.Cells.Item(132,1).Value  = 10
$ws.Cells.Item(132,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(132,3).Value  = "La Araucanía"
$ws.Cells.Item(132,4).Value  = 44588
$ws.Cells.Item(132,5).Value  = 9
$ws.Cells.Item(132,6).Value  = "Fruta"
$ws.Cells.Item(132,7).Value  = 100103
$ws.Cells.Item(132,8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(132,9).Value  = 100103002
$ws.Cells.Item(132,10).Value = "Ciruela"
$ws.Cells.Item(132,11).Value = "Black Amber"
$ws.Cells.Item(132,12).Value = "Primera"
$ws.Cells.Item(132,13).Value = 330
$ws.Cells.Item(132,14).Value = 14000
$ws.Cells.Item(132,15).Value = 15000
$ws.Cells.Item(132,16).Value = 14545
$ws.Cells.Item(132,17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(132,18).Value = "Región de O'Higgins"
$ws.Cells.Item(132,19).Value = 808
$ws.Cells.Item(132,20).Value = 18
